# Update stats for players (rows 2-11 on Sheet1) with refreshed session data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: D = NoOfSessions, I = PointsBonusTotal, J = AvgSessPoints

$ws.Range("D2").Value = 127
$ws.Range("I2").Value = 460
$ws.Range("J2").Value = 3.62

$ws.Range("D3").Value = 147
$ws.Range("I3").Value = 525
$ws.Range("J3").Value = 3.57

$ws.Range("D4").Value = 197
$ws.Range("I4").Value = 736
$ws.Range("J4").Value = 3.74

$ws.Range("D5").Value = 109
$ws.Range("I5").Value = 398
$ws.Range("J5").Value = 3.65

$ws.Range("D6").Value = 44
$ws.Range("I6").Value = 161
$ws.Range("J6").Value = 3.66

$ws.Range("D7").Value = 148
$ws.Range("I7").Value = 632

$ws.Range("D8").Value = 202
$ws.Range("I8").Value = 705
$ws.Range("J8").Value = 3.49

$ws.Range("D9").Value = 88
$ws.Range("I9").Value = 300
$ws.Range("J9").Value = 3.41

$ws.Range("D10").Value = 205
$ws.Range("I10").Value = 847
$ws.Range("J10").Value = 4.13

$ws.Range("D11").Value = 124
$ws.Range("I11").Value = 487
